$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.156.60"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "3.809.15"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "702.42"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.51"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").Value = "3.808.80"
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.52"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.88"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").Value = "4.452.90"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "3.825.87"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").Value = "71.164.36"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.114"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.14"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "513.27"
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("E25").Value = "  -3.35%  "
$ws.Range("D26").Value = "3.958.77"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.06"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.35"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -3.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.02"
$ws.Range("E31").Value = "  -5.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.39"
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.01"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").Value = "3.772.76"
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.994"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.38"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.28"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "174.00"
$ws.Range("E45").Value = "  +6.49%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000311"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.34"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "426.61"
$ws.Range("E49").Value = "  +3.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.37"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.60"
$ws.Range("E51").Value = "  -0.20%  "
